# Reorganized and knitted scripts
#
# The workbook originally has 20 sheets named "250m", "500m", ... "5000m",
# each holding the same covariate-correlation template (a header row of
# term names plus the term names repeated down column A).
#
# This edit inserts two additional buffer-distance sheets, "50m" and
# "100m", at the very front of the tab order (before "250m"), each
# carrying the identical template used by every other sheet. All existing
# sheets and their data are left untouched.

$wb = $excel.ActiveWorkbook

# The first sheet ("250m") holds the template content/style that every
# sheet in this workbook shares - duplicate it to seed the two new sheets
# so the header formatting (bold + centered) and term list come along for
# free, then just rename/position the copies.

$firstSheet = $wb.Worksheets.Item(1)

# Duplicate "250m" and drop the copy immediately before it -> becomes the
# new first tab; rename it to "100m".
$firstSheet.Copy($firstSheet)
$sheet100m = $wb.Worksheets.Item(1)
$sheet100m.Name = "100m"

# Duplicate the freshly-made "100m" and drop that copy before it -> new
# first tab; rename it to "50m". Final order: 50m, 100m, 250m, 500m, ...
$sheet100m.Copy($sheet100m)
$sheet50m = $wb.Worksheets.Item(1)
$sheet50m.Name = "50m"

# Keep the first tab active/selected, matching a freshly reorganized
# workbook opening on its first sheet.
$sheet50m.Select()
